$d = $word.ActiveDocument

# --- 1. Title ---
$null = $d.Content.Find.Execute(
    "Harmony of Numbers: Unveiling Mathematical Patterns in Nature",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "The Alluring Symphony of Chemistry: Unveiling the Essence of Matter", 2)

# --- 2. Author name ---
$null = $d.Content.Find.Execute(
    "Sophia Carter",
    $false, $true, $false, $false, $false, $true, 1, $false,
    "Emily Davis", 2)

# --- 3. Email paragraph (sophiacarter@mathland / . / institute -> EmilyDavis@schoolmail / . / org) ---
# Replaced wholesale via raw OOXML so the run boundaries/count stay exactly as authored.
$xmlPara3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BB7161" w:rsidRDefault="00646610"><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>EmilyDavis@schoolmail</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="32"/></w:rPr><w:t>org</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(3).Range.InsertXML($xmlPara3)

# --- 5. Big body paragraph: many sentence-level rewrites + two 3-run -> 1-run merges ---
$xmlPara5 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BB7161" w:rsidRDefault="00646610"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>Chemistry, a science that delves into the intricacies of matter and its transformations, offers a captivating journey into the microscopic realm of atoms and molecules</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> This fascinating field unveils the secrets of chemical reactions, enabling us to comprehend the composition of substances and the mechanisms behind their interactions</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Chemistry permeates our daily lives, from the food we consume to the medicines we rely on, shaping our world in profound ways</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>As we embark on this exploration of chemistry, we will discover the fundamental principles that govern the behavior ofWu Zhi , including the periodic table, atomic structure, and bonding</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> We will delve into the captivating world of chemical reactions, unraveling the mysteries of how substances interact and transform into new substances</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Furthermore, we will explore the practical applications of chemistry in various fields, witnessing its indispensable role in industries such as pharmaceuticals, materials science, and energy</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:br/><w:t>The study of chemistry is not merely an accumulation of knowledge; it is an invitation to embark on an intellectual adventure, fostering critical thinking and problem-solving skills</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> By unraveling the intricacies of chemical phenomena, we cultivate a deeper understanding of our physical world and the processes that shape it</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> Chemistry empowers us to address global challenges, such as developing sustainable energy sources and combating environmental pollution, making it a discipline of paramount importance in our contemporary world</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(5).Range.InsertXML($xmlPara5)

# --- 7. Summary paragraph: rewrites + dropped lastRenderedPageBreak + new run split ---
$xmlPara7 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00BB7161" w:rsidRDefault="00646610"><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>Chemistry, a captivating science that unravels the secrets of matter and its transformations, plays a pivotal role in shaping our world</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Through the study of chemistry, we gain a deeper understanding of the composition, behavior, and interactions of substances</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> This knowledge has wide-ranging applications, spanning industries and </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:lastRenderedPageBreak/><w:t>impacting our daily lives</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t xml:space="preserve"> Chemistry empowers us to address global challenges and fosters critical thinking skills, making it an essential discipline in the 21st century</w:t></w:r><w:r w:rsidR="006A51D7"><w:rPr><w:rFonts w:ascii="Aptos" w:hAnsi="Aptos"/><w:color w:val="000000"/></w:rPr><w:t>.</w:t></w:r></w:p>
<w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$d.Paragraphs(7).Range.InsertXML($xmlPara7)
